# Applies weekly re-shuffle of Fruta/Hortaliza records (Chirimoya) rows 2-10
# Column layout: A..T -> D=Fecha(col4) L=Calidad(col12) M=Volumen(col13)
# N=Precio minimo(col14) O=Precio maximo(col15) P=Precio promedio(col16)
# Q=Unidad comercializacion(col17) S=Precio $/Kg(col19) T=Kg/unidad(col20)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = @{ D = 44475; L = "Especial"; M = 200; N = 32000; O = 33000; P = 32500; Q = "`$/caja 12 kilos"; S = 2708; T = 12 }
    3  = @{ D = 44489; L = "Primera";  M = 200; N = 24000; O = 25000; P = 24500; Q = "`$/caja 12 kilos"; S = 2042; T = 12 }
    4  = @{ D = 44482; L = "Primera";  M = 160; N = 25000; O = 26000; P = 25500; Q = "`$/caja 12 kilos"; S = 2125; T = 12 }
    5  = @{ D = 44167; L = "Segunda";  M = 200; N = 18000; O = 19000; P = 18500; Q = "`$/caja 13 kilos"; S = 1423; T = 13 }
    6  = @{ D = 44441; L = "Primera";  M = 100; N = 29000; O = 30000; P = 29500; Q = "`$/caja 12 kilos"; S = 2458; T = 12 }
    8  = @{ D = 44524; L = "Primera";  M = 200; N = 23000; O = 24000; P = 23500; Q = "`$/caja 12 kilos"; S = 1958; T = 12 }
    9  = @{ D = 44545; L = "Primera";  M = 200; N = 23000; O = 24000; P = 23500; Q = "`$/bandeja 12 kilos"; S = 1958; T = 12 }
    10 = @{ D = 44496; L = "Primera";  M = 200; N = 23000; O = 24000; P = 23500; Q = "`$/caja 12 kilos"; S = 1958; T = 12 }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 4).Value  = $vals.D   # D - Fecha
    $ws.Cells.Item($row, 12).Value = $vals.L   # L - Calidad
    $ws.Cells.Item($row, 13).Value = $vals.M   # M - Volumen
    $ws.Cells.Item($row, 14).Value = $vals.N   # N - Precio minimo
    $ws.Cells.Item($row, 15).Value = $vals.O   # O - Precio maximo
    $ws.Cells.Item($row, 16).Value = $vals.P   # P - Precio promedio ponderado
    $ws.Cells.Item($row, 17).Value = $vals.Q   # Q - Unidad de comercializacion
    $ws.Cells.Item($row, 19).Value = $vals.S   # S - Precio $/Kg
    $ws.Cells.Item($row, 20).Value = $vals.T   # T - Kg / unidad
}

$wb.Save()
